$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (existing header cell) onto the two new
# header cells so they pick up the same cellXf (bold font, border,
# centered/top alignment) instead of Excel generating a brand new style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header values
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows (plain numeric cells, no style like the other data columns)
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 8

$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 6
